$wb = $excel.ActiveWorkbook

# Grab the existing "Message" sheet up front - new sheets are always added
# after the current last sheet (below), so this reference stays valid for
# the whole script.
$msg = $wb.Worksheets.Item("Message")

# ---------------------------------------------------------------------------
# New Cypher query used by the "stat" sheet (counts of files/samples/cases/
# studies) added alongside the existing filtered-record query.
# ---------------------------------------------------------------------------
$statCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN ['Female']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# ---------------------------------------------------------------------------
# A scratch sheet (added & removed before the permanent sheets, so the
# sheetId/relationship numbering of the final sheets is not disturbed) is
# used to stage the numeric-looking values ("214", "74", ...) as
# quote-prefixed text. A values-only paste from it then writes those values
# into the real sheet as plain text (shared string) cells rather than
# numbers, without leaving any cell-level number formatting behind on the
# destination cells.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch = $wb.Worksheets.Add($null, $lastSheet)
$scratch.Cells.Item(1, 1).Value = "'214"
$scratch.Cells.Item(1, 2).Value = "'74"
$scratch.Cells.Item(1, 3).Value = "'33"
$scratch.Cells.Item(1, 4).Value = "'1"
$scratch.Range("A1:D1").Copy()
$scratch.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: CypherOutput_Message - an exact duplicate of the "Message" sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"
for ($r = 1; $r -le 10; $r++) {
    $cypherOutputMessage.Cells.Item($r, 1).Value = $msg.Cells.Item($r, 1).Text
}

# ---------------------------------------------------------------------------
# Sheet 4: StatOutput - header row of metric names + one row of counts.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"
$statOutput.Cells.Item(1, 1).Value = "number_of_files"
$statOutput.Cells.Item(1, 2).Value = "number_of_sample"
$statOutput.Cells.Item(1, 3).Value = "number_of_cases"
$statOutput.Cells.Item(1, 4).Value = "number_of_study"
$statOutput.Range("A2").PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# Sheet 5: StatOutput_Message - same connection-info block as "Message",
# repeated twice; the second copy's Cypher line (row 18) is the new stat
# query instead of the original filter query.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$statOutputMessage.Name = "StatOutput_Message"
for ($r = 1; $r -le 10; $r++) {
    $statOutputMessage.Cells.Item($r, 1).Value = $msg.Cells.Item($r, 1).Text
}
for ($r = 1; $r -le 7; $r++) {
    $statOutputMessage.Cells.Item(10 + $r, 1).Value = $msg.Cells.Item($r, 1).Text
}
$statOutputMessage.Cells.Item(18, 1).Value = $statCypher
$statOutputMessage.Cells.Item(19, 1).Value = $msg.Cells.Item(9, 1).Text
$statOutputMessage.Cells.Item(20, 1).Value = $msg.Cells.Item(10, 1).Text
